$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45174 (2023-09-05) to 45175 (2023-09-06) for rows 2 through 19.
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45175
}
